$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'29.435.36"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "'1.868.43"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'243.63"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").Value = "'0.7044"
$ws.Range("E6").Value = "  -1.11%  "

$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.3143"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("D9").Value = "'0.07851"
$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").Value = "'24.49"
$ws.Range("E10").Value = "  -2.42%  "

$ws.Range("D11").Value = "'0.08020"
$ws.Range("E11").Value = "  -3.79%  "

$ws.Range("D12").Value = "'1.893.50"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "'5.198"
$ws.Range("E13").Value = "  -1.37%  "

$ws.Range("D14").Value = "'93.54"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("D15").Value = "'0.7014"
$ws.Range("E15").Value = "  -2.41%  "

$ws.Range("D16").Value = "'6.463"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").Value = "'29.522.54"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "'0.000008319"
$ws.Range("E18").Value = "  -4.10%  "

$ws.Range("D19").Value = "'255.22"
$ws.Range("E19").Value = "  +5.07%  "

$ws.Range("D20").Value = "'2.139.97"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").Value = "'13.15"
$ws.Range("E21").Value = "  -1.48%  "

$ws.Range("D23").Value = "'7.597"
$ws.Range("E23").Value = "  -3.22%  "

$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").Value = "'0.1555"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("D26").Value = "'9.040"
$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("D27").Value = "'161.07"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("D29").Value = "'1.501"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").Value = "'4.324"
$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("D31").Value = "'4.266"
$ws.Range("E31").Value = "  -2.02%  "

$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").Value = "'0.05314"
$ws.Range("E33").Value = "  -1.46%  "

$ws.Range("D34").Value = "'1.890"
$ws.Range("E34").Value = "  -3.03%  "

$ws.Range("D35").Value = "'0.7460"
$ws.Range("E35").Value = "  -3.76%  "

$ws.Range("D36").Value = "'1.164"
$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("E37").Value = "  +1.28%  "

$ws.Range("D38").Value = "'0.01875"
$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("D39").Value = "'1.260.44"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").Value = "'2.743"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").Value = "'0.8987"
$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("D42").Value = "'108.82"
$ws.Range("E42").Value = "  -3.85%  "

$ws.Range("D43").Value = "'5.947"
$ws.Range("E43").Value = "  -8.86%  "

$ws.Range("D44").Value = "'71.32"
$ws.Range("E44").Value = "  -4.34%  "

$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "'0.00000000130"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("D47").Value = "'2.038.00"

$ws.Range("D48").Value = "'1.798"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("D49").Value = "'0.5190"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").Value = "'9.487"
$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("D51").Value = "'0.4308"
$ws.Range("E51").Value = "  -1.71%  "
